$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray "x" / helper-note cells that are no longer needed.
$ws.Range("E22").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("D30").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("C61").ClearContents()
$ws.Range("D61").ClearContents()

# The cells that lost their "x" marker had their font color reset from red
# back to the automatic/default color.
$resetFontCells = @("D22", "B24", "B25", "B26", "B61")
foreach ($cellRef in $resetFontCells) {
    $ws.Range($cellRef).Font.ThemeColor = 1
}

# Update the saved view state: scroll/selection moved to F22.
$ws.Range("F22").Select()
